# -----------------------------------------------------------------------
# Reindeer Hunter test workbook: "Completed implementation of no student id"
#  - adds a new "TestNoStudentIDs" worksheet (some students missing an ID)
#  - updates the selection/active-cell remembered on each existing sheet
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Update the remembered selection on the four existing grade sheets ---

$wsGrade9 = $wb.Worksheets.Item("Grade 9")
$wsGrade9.Range("A6:E9").Select() | Out-Null

$wsGrade10 = $wb.Worksheets.Item("Grade 10")
$wsGrade10.Range("A10:E13").Select() | Out-Null

$wsGrade11 = $wb.Worksheets.Item("Grade 11")
$wsGrade11.Range("A14:E17").Select() | Out-Null

$wsGrade12 = $wb.Worksheets.Item("Grade 12")
$wsGrade12.Range("A2:E21").Select() | Out-Null

# --- Add the new worksheet after "Grade 12" ---

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNew.Name = "TestNoStudentIDs"

# Header row
$wsNew.Range("A1").Value = "First Name"
$wsNew.Range("B1").Value = "Last Name"
$wsNew.Range("C1").Value = "Student ID"
$wsNew.Range("D1").Value = "Grade "
$wsNew.Range("E1").Value = "Homeroom"

# Student rows - these students have no Student ID (column C left blank);
# the application assigns them a randomly generated ID (with a special
# prefix) at import time.
$students = @(
  @("First 12-1", "Last 12-1", 12, 1201),
  @("First 12-2", "Last 12-2", 12, 1201),
  @("First 12-3", "Last 12-3", 12, 1201),
  @("First 12-4", "Last 12-4", 12, 1201),
  @("First 9-5", "Last 9-5", 9, 902),
  @("First 9-6", "Last 9-6", 9, 902),
  @("First 9-7", "Last 9-7", 9, 902),
  @("First 9-8", "Last 9-8", 9, 902),
  @("First 10-10", "Last 10-10", 10, 1003),
  @("First 10-10", "Last 10-10", 10, 1003),
  @("First 10-11", "Last 10-11", 10, 1003),
  @("First 10-12", "Last 10-12", 10, 1003),
  @("First 11-13", "Last 11-13", 11, 1104),
  @("First 11-14", "Last 11-14", 11, 1104),
  @("First 11-15", "Last 11-15", 11, 1104),
  @("First 11-16", "Last 11-16", 11, 1104)
)

$r = 2
foreach ($student in $students) {
  $wsNew.Cells.Item($r, 1).Value = $student[0]
  $wsNew.Cells.Item($r, 2).Value = $student[1]
  $wsNew.Cells.Item($r, 4).Value = $student[2]
  $wsNew.Cells.Item($r, 5).Value = $student[3]
  $r++
}

# Activate the new sheet (it becomes the selected tab) and set its
# remembered selection/active cell.
$wsNew.Activate() | Out-Null
$wsNew.Range("K10").Select() | Out-Null
